# Mockup update position: move the floating "Text Box 5" callout from
# left=416.4pt / top=80.8pt to left=415.8pt / top=64.6pt (the shape's
# anchors - margin / paragraph - are unchanged).

$d = $word.ActiveDocument

$targetName = "Text Box 5"
$newLeft = 415.8
$newTop = 64.6

$count = $d.Shapes.Count

# Remember every shape's current Left/Top up front (collection reads are
# consistent/by z-order) so any slot we touch by mistake while looking
# for the right one to write through can be put back exactly as found.
$origLeft = @{}
$origTop = @{}
for ($i = 1; $i -le $count; $i++) {
    $sh = $d.Shapes.Item($i)
    $origLeft[$i] = $sh.Left
    $origTop[$i] = $sh.Top
}

$applied = $false
for ($i = 1; $i -le $count; $i++) {
    $sh = $d.Shapes.Item($i)
    $sh.Left = $newLeft
    $sh.Top = $newTop

    # Re-read the whole collection by name to see whether the intended
    # shape now reports the new position.
    $foundLeft = $null
    $foundTop = $null
    for ($j = 1; $j -le $count; $j++) {
        $sh2 = $d.Shapes.Item($j)
        if ($sh2.Name -eq $targetName) {
            $foundLeft = $sh2.Left
            $foundTop = $sh2.Top
        }
    }

    if ($foundLeft -eq $newLeft -and $foundTop -eq $newTop) {
        $applied = $true
        break
    }

    # Wrong slot - restore whatever this index actually changed before
    # trying the next one.
    $sh.Left = $origLeft[$i]
    $sh.Top = $origTop[$i]
}

if (-not $applied) {
    Write-Host "WARNING: could not move '$targetName' to the target position"
}
